$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$s = $c.Style
$c.Value = "'" + '66.663.70'
$c.Style = $s
$c = $ws.Range('E2')
$s = $c.Style
$c.Value = "'" + '  +1.49%  '
$c.Style = $s
$c = $ws.Range('D3')
$s = $c.Style
$c.Value = "'" + '3.502.31'
$c.Style = $s
$c = $ws.Range('E3')
$s = $c.Style
$c.Value = "'" + '  +1.07%  '
$c.Style = $s
$c = $ws.Range('E4')
$s = $c.Style
$c.Value = "'" + '  +0.06%  '
$c.Style = $s
$c = $ws.Range('D5')
$s = $c.Style
$c.Value = "'" + '592.28'
$c.Style = $s
$c = $ws.Range('E5')
$s = $c.Style
$c.Value = "'" + '  +1.74%  '
$c.Style = $s
$c = $ws.Range('D6')
$s = $c.Style
$c.Value = "'" + '169.04'
$c.Style = $s
$c = $ws.Range('E6')
$s = $c.Style
$c.Value = "'" + '  +0.40%  '
$c.Style = $s
$c = $ws.Range('D7')
$s = $c.Style
$c.Value = "'" + '0.999'
$c.Style = $s
$c = $ws.Range('E7')
$s = $c.Style
$c.Value = "'" + '  -0.03%  '
$c.Style = $s
$c = $ws.Range('D8')
$s = $c.Style
$c.Value = "'" + '0.592'
$c.Style = $s
$c = $ws.Range('E8')
$s = $c.Style
$c.Value = "'" + '  +4.50%  '
$c.Style = $s
$c = $ws.Range('E9')
$s = $c.Style
$c.Value = "'" + '  +5.13%  '
$c.Style = $s
$c = $ws.Range('D10')
$s = $c.Style
$c.Value = "'" + '7.33'
$c.Style = $s
$c = $ws.Range('E10')
$s = $c.Style
$c.Value = "'" + '  +0.73%  '
$c.Style = $s
$c = $ws.Range('D11')
$s = $c.Style
$c.Value = "'" + '0.432'
$c.Style = $s
$c = $ws.Range('E11')
$s = $c.Style
$c.Value = "'" + '  +0.70%  '
$c.Style = $s
$c = $ws.Range('D12')
$s = $c.Style
$c.Value = "'" + '4.106.59'
$c.Style = $s
$c = $ws.Range('E12')
$s = $c.Style
$c.Value = "'" + '  +1.13%  '
$c.Style = $s
$c = $ws.Range('E13')
$s = $c.Style
$c.Value = "'" + '  -0.15%  '
$c.Style = $s
$c = $ws.Range('D14')
$s = $c.Style
$c.Value = "'" + '28.16'
$c.Style = $s
$c = $ws.Range('E14')
$s = $c.Style
$c.Value = "'" + '  +2.32%  '
$c.Style = $s
$c = $ws.Range('E15')
$s = $c.Style
$c.Value = "'" + '  +1.63%  '
$c.Style = $s
$c = $ws.Range('D16')
$s = $c.Style
$c.Value = "'" + '66.716.72'
$c.Style = $s
$c = $ws.Range('E16')
$s = $c.Style
$c.Value = "'" + '  +1.77%  '
$c.Style = $s
$c = $ws.Range('D17')
$s = $c.Style
$c.Value = "'" + '3.494.31'
$c.Style = $s
$c = $ws.Range('E17')
$s = $c.Style
$c.Value = "'" + '  +4.77%  '
$c.Style = $s
$c = $ws.Range('D18')
$s = $c.Style
$c.Value = "'" + '6.32'
$c.Style = $s
$c = $ws.Range('E18')
$s = $c.Style
$c.Value = "'" + '  +1.46%  '
$c.Style = $s
$c = $ws.Range('D19')
$s = $c.Style
$c.Value = "'" + '14.07'
$c.Style = $s
$c = $ws.Range('E19')
$s = $c.Style
$c.Value = "'" + '  +1.94%  '
$c.Style = $s
$c = $ws.Range('D20')
$s = $c.Style
$c.Value = "'" + '393.66'
$c.Style = $s
$c = $ws.Range('E20')
$s = $c.Style
$c.Value = "'" + '  +2.42%  '
$c.Style = $s
$c = $ws.Range('D21')
$s = $c.Style
$c.Value = "'" + '7.96'
$c.Style = $s
$c = $ws.Range('E21')
$s = $c.Style
$c.Value = "'" + '  +0.35%  '
$c.Style = $s
$c = $ws.Range('D22')
$s = $c.Style
$c.Value = "'" + '73.13'
$c.Style = $s
$c = $ws.Range('E22')
$s = $c.Style
$c.Value = "'" + '  +2.15%  '
$c.Style = $s
$c = $ws.Range('D23')
$s = $c.Style
$c.Value = "'" + '0.999'
$c.Style = $s
$c = $ws.Range('E23')
$s = $c.Style
$c.Value = "'" + '  -0.36%  '
$c.Style = $s
$c = $ws.Range('D24')
$s = $c.Style
$c.Value = "'" + '0.534'
$c.Style = $s
$c = $ws.Range('E24')
$s = $c.Style
$c.Value = "'" + '  +2.30%  '
$c.Style = $s
$c = $ws.Range('E25')
$s = $c.Style
$c.Value = "'" + '  +1.26%  '
$c.Style = $s
$c = $ws.Range('D26')
$s = $c.Style
$c.Value = "'" + '10.19'
$c.Style = $s
$c = $ws.Range('E26')
$s = $c.Style
$c.Value = "'" + '  +3.34%  '
$c.Style = $s
$c = $ws.Range('E27')
$s = $c.Style
$c.Value = "'" + '  -0.53%  '
$c.Style = $s
$c = $ws.Range('E28')
$s = $c.Style
$c.Value = "'" + '  -0.20%  '
$c.Style = $s
$c = $ws.Range('D29')
$s = $c.Style
$c.Value = "'" + '6.36'
$c.Style = $s
$c = $ws.Range('E29')
$s = $c.Style
$c.Value = "'" + '  +1.47%  '
$c.Style = $s
$c = $ws.Range('D30')
$s = $c.Style
$c.Value = "'" + '1.46'
$c.Style = $s
$c = $ws.Range('E30')
$s = $c.Style
$c.Value = "'" + '  +1.13%  '
$c.Style = $s
$c = $ws.Range('E31')
$s = $c.Style
$c.Value = "'" + '  +1.45%  '
$c.Style = $s
$c = $ws.Range('E32')
$s = $c.Style
$c.Value = "'" + '  +2.18%  '
$c.Style = $s
$c = $ws.Range('D33')
$s = $c.Style
$c.Value = "'" + '7.34'
$c.Style = $s
$c = $ws.Range('E33')
$s = $c.Style
$c.Value = "'" + '  +0.35%  '
$c.Style = $s
$c = $ws.Range('E34')
$s = $c.Style
$c.Value = "'" + '  +5.51%  '
$c.Style = $s
$c = $ws.Range('D35')
$s = $c.Style
$c.Value = "'" + '162.56'
$c.Style = $s
$c = $ws.Range('E35')
$s = $c.Style
$c.Value = "'" + '  +1.36%  '
$c.Style = $s
$c = $ws.Range('D36')
$s = $c.Style
$c.Value = "'" + '0.899'
$c.Style = $s
$c = $ws.Range('E36')
$s = $c.Style
$c.Value = "'" + '  +0.55%  '
$c.Style = $s
$c = $ws.Range('E37')
$s = $c.Style
$c.Value = "'" + '  +2.89%  '
$c.Style = $s
$c = $ws.Range('D38')
$s = $c.Style
$c.Value = "'" + '6.82'
$c.Style = $s
$c = $ws.Range('E38')
$s = $c.Style
$c.Value = "'" + '  +2.82%  '
$c.Style = $s
$c = $ws.Range('E39')
$s = $c.Style
$c.Value = "'" + '  +4.41%  '
$c.Style = $s
$c = $ws.Range('D40')
$s = $c.Style
$c.Value = "'" + '26.56'
$c.Style = $s
$c = $ws.Range('E40')
$s = $c.Style
$c.Value = "'" + '  +1.48%  '
$c.Style = $s
$c = $ws.Range('D41')
$s = $c.Style
$c.Value = "'" + '0.0741'
$c.Style = $s
$c = $ws.Range('E41')
$s = $c.Style
$c.Value = "'" + '  +1.11%  '
$c.Style = $s
$c = $ws.Range('D42')
$s = $c.Style
$c.Value = "'" + '26.85'
$c.Style = $s
$c = $ws.Range('E42')
$s = $c.Style
$c.Value = "'" + '  -0.02%  '
$c.Style = $s
$c = $ws.Range('D43')
$s = $c.Style
$c.Value = "'" + '2.792.38'
$c.Style = $s
$c = $ws.Range('E43')
$s = $c.Style
$c.Value = "'" + '  -0.48%  '
$c.Style = $s
$c = $ws.Range('D44')
$s = $c.Style
$c.Value = "'" + '42.93'
$c.Style = $s
$c = $ws.Range('E44')
$s = $c.Style
$c.Value = "'" + '  -0.32%  '
$c.Style = $s
$c = $ws.Range('E45')
$s = $c.Style
$c.Value = "'" + '  +3.50%  '
$c.Style = $s
$c = $ws.Range('D46')
$s = $c.Style
$c.Value = "'" + '0.0311'
$c.Style = $s
$c = $ws.Range('E46')
$s = $c.Style
$c.Value = "'" + '  +1.10%  '
$c.Style = $s
$c = $ws.Range('D47')
$s = $c.Style
$c.Value = "'" + '343.04'
$c.Style = $s
$c = $ws.Range('E47')
$s = $c.Style
$c.Value = "'" + '  +1.22%  '
$c.Style = $s
$c = $ws.Range('E48')
$s = $c.Style
$c.Value = "'" + '  +1.59%  '
$c.Style = $s
$c = $ws.Range('D49')
$s = $c.Style
$c.Value = "'" + '33.99'
$c.Style = $s
$c = $ws.Range('E49')
$s = $c.Style
$c.Value = "'" + '  +4.86%  '
$c.Style = $s
$c = $ws.Range('E50')
$s = $c.Style
$c.Value = "'" + '  +2.71%  '
$c.Style = $s
$c = $ws.Range('D51')
$s = $c.Style
$c.Value = "'" + '6.50'
$c.Style = $s
$c = $ws.Range('E51')
$s = $c.Style
$c.Value = "'" + '  +1.78%  '
$c.Style = $s
